$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap differing cells between rows 2 and 3
$ws.Range("A2").Value2 = 131046824
$ws.Range("A3").Value2 = 131046825
$ws.Range("Q2").Value2 = 401653
$ws.Range("Q3").Value2 = 401650
$ws.Range("R2").Value2 = 6818054
$ws.Range("R3").Value2 = 6818017
$ws.Range("Z2").Value2 = "14:50"
$ws.Range("Z3").Value2 = "14:52"
$ws.Range("AB2").Value2 = "14:50"
$ws.Range("AB3").Value2 = "14:52"

# Swap differing cells between rows 10 and 11
$ws.Range("A10").Value2 = 131046823
$ws.Range("A11").Value2 = 131046773
$ws.Range("B10").Value2 = 79243
$ws.Range("B11").Value2 = 57884
$ws.Range("E10").Value2 = 6425
$ws.Range("E11").Value2 = 100109
$ws.Range("F10").Value2 = "Garnlav"
$ws.Range("F11").Value2 = "Tretåig hackspett"
$ws.Range("G10").Value2 = "Alectoria sarmentosa"
$ws.Range("G11").Value2 = "Picoides tridactylus"
$ws.Range("H10").Value2 = "(Ach.) Ach."
$ws.Range("H11").Value2 = "(Linnaeus, 1758)"
$ws.Range("M10").Value2 = ""
$ws.Range("M11").Value2 = "äldre spår"
$ws.Range("Q10").Value2 = 401661
$ws.Range("Q11").Value2 = 401346
$ws.Range("R10").Value2 = 6818064
$ws.Range("R11").Value2 = 6818162
$ws.Range("Z10").Value2 = "14:50"
$ws.Range("Z11").Value2 = "15:23"
$ws.Range("AB10").Value2 = "14:50"
$ws.Range("AB11").Value2 = "15:23"
$ws.Range("AC10").Value2 = ""
$ws.Range("AC11").Value2 = "Äldre ringhack (gran)"

# Swap differing cells between rows 12 and 13
$ws.Range("A12").Value2 = 131046828
$ws.Range("A13").Value2 = 131046769
$ws.Range("B12").Value2 = 79243
$ws.Range("B13").Value2 = 57884
$ws.Range("E12").Value2 = 6425
$ws.Range("E13").Value2 = 100109
$ws.Range("F12").Value2 = "Garnlav"
$ws.Range("F13").Value2 = "Tretåig hackspett"
$ws.Range("G12").Value2 = "Alectoria sarmentosa"
$ws.Range("G13").Value2 = "Picoides tridactylus"
$ws.Range("H12").Value2 = "(Ach.) Ach."
$ws.Range("H13").Value2 = "(Linnaeus, 1758)"
$ws.Range("M12").Value2 = ""
$ws.Range("M13").Value2 = "färska spår"
$ws.Range("Q12").Value2 = 401634
$ws.Range("Q13").Value2 = 401575
$ws.Range("R12").Value2 = 6817871
$ws.Range("R13").Value2 = 6817873
$ws.Range("Z12").Value2 = "14:58"
$ws.Range("Z13").Value2 = "15:07"
$ws.Range("AB12").Value2 = "14:58"
$ws.Range("AB13").Value2 = "15:07"
$ws.Range("AC12").Value2 = ""
$ws.Range("AC13").Value2 = "Färska ringhack (tall)"

# Swap differing cells between rows 14 and 15
$ws.Range("A14").Value2 = 131046771
$ws.Range("A15").Value2 = 131046709
$ws.Range("B14").Value2 = 57884
$ws.Range("B15").Value2 = 83223
$ws.Range("E14").Value2 = 100109
$ws.Range("E15").Value2 = 6440
$ws.Range("F14").Value2 = "Tretåig hackspett"
$ws.Range("F15").Value2 = "Vitgrynig nållav"
$ws.Range("G14").Value2 = "Picoides tridactylus"
$ws.Range("G15").Value2 = "Chaenotheca subroscida"
$ws.Range("H14").Value2 = "(Linnaeus, 1758)"
$ws.Range("H15").Value2 = "(Eitner) Zahlbr."
$ws.Range("M14").Value2 = "färska spår"
$ws.Range("M15").Value2 = ""
$ws.Range("Q14").Value2 = 401556
$ws.Range("Q15").Value2 = 401646
$ws.Range("R14").Value2 = 6817954
$ws.Range("R15").Value2 = 6817967
$ws.Range("Z14").Value2 = "15:12"
$ws.Range("Z15").Value2 = "14:54"
$ws.Range("AB14").Value2 = "15:12"
$ws.Range("AB15").Value2 = "14:54"
$ws.Range("AC14").Value2 = "Färska ringhack (tall)"
$ws.Range("AC15").Value2 = ""

# Swap differing cells between rows 20 and 21
$ws.Range("A20").Value2 = 131047034
$ws.Range("A21").Value2 = 131046766
$ws.Range("B20").Value2 = 78646
$ws.Range("B21").Value2 = 57884
$ws.Range("E20").Value2 = 6437
$ws.Range("E21").Value2 = 100109
$ws.Range("F20").Value2 = "Blanksvart spiklav"
$ws.Range("F21").Value2 = "Tretåig hackspett"
$ws.Range("G20").Value2 = "Calicium denigratum"
$ws.Range("G21").Value2 = "Picoides tridactylus"
$ws.Range("H20").Value2 = "(Vain.) Tibell"
$ws.Range("H21").Value2 = "(Linnaeus, 1758)"
$ws.Range("K20").Value2 = ""
$ws.Range("K21").Value2 = ""
$ws.Range("L20").Value2 = ""
$ws.Range("L21").Value2 = ""
$ws.Range("M20").Value2 = ""
$ws.Range("M21").Value2 = "äldre spår"
$ws.Range("N20").Value2 = ""
$ws.Range("N21").Value2 = ""
$ws.Range("Q20").Value2 = 401597
$ws.Range("Q21").Value2 = 401318
$ws.Range("R20").Value2 = 6817852
$ws.Range("R21").Value2 = 6818379
$ws.Range("Z20").Value2 = "15:05"
$ws.Range("Z21").Value2 = "15:29"
$ws.Range("AB20").Value2 = "15:05"
$ws.Range("AB21").Value2 = "15:29"
$ws.Range("AC20").Value2 = ""
$ws.Range("AC21").Value2 = "Äldre ringhack (tall)"

# Swap differing cells between rows 26 and 27
$ws.Range("A26").Value2 = 131046832
$ws.Range("A27").Value2 = 131047014
$ws.Range("B26").Value2 = 79243
$ws.Range("B27").Value2 = 57884
$ws.Range("E26").Value2 = 6425
$ws.Range("E27").Value2 = 100109
$ws.Range("F26").Value2 = "Garnlav"
$ws.Range("F27").Value2 = "Tretåig hackspett"
$ws.Range("G26").Value2 = "Alectoria sarmentosa"
$ws.Range("G27").Value2 = "Picoides tridactylus"
$ws.Range("H26").Value2 = "(Ach.) Ach."
$ws.Range("H27").Value2 = "(Linnaeus, 1758)"
$ws.Range("M26").Value2 = ""
$ws.Range("M27").Value2 = "färska spår"
$ws.Range("Q26").Value2 = 401350
$ws.Range("Q27").Value2 = 401378
$ws.Range("R26").Value2 = 6818162
$ws.Range("R27").Value2 = 6818082
$ws.Range("Z26").Value2 = "15:24"
$ws.Range("Z27").Value2 = "15:21"
$ws.Range("AB26").Value2 = "15:24"
$ws.Range("AB27").Value2 = "15:21"
$ws.Range("AC26").Value2 = ""
$ws.Range("AC27").Value2 = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE26").Value2 = $False
$ws.Range("AE27").Value2 = $True
